$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 5.683658814236868
$ws.Range("D2").Value = 8.279395454628119
$ws.Range("E2").Value = 13.30210991916849
$ws.Range("F2").Value = 37.11039070033118
$ws.Range("G2").Value = 42.02745061849008
$ws.Range("H2").Value = 17.40151863775981
$ws.Range("I2").Value = 20.36751465146092
$ws.Range("J2").Value = 10.07129936889861
$ws.Range("K2").Value = 15.10813732258622
$ws.Range("N2").Value = 18.87054176394461

$ws.Range("B3").Value = 5.574485982242754
$ws.Range("D3").Value = 8.246043383262229
$ws.Range("E3").Value = 13.25377573637131
$ws.Range("F3").Value = 37.06822747821502
$ws.Range("G3").Value = 41.90373264459933
$ws.Range("H3").Value = 17.4360548641228
$ws.Range("I3").Value = 20.46446283322906
$ws.Range("J3").Value = 10.07863977131413
$ws.Range("K3").Value = 14.77623572043987
$ws.Range("N3").Value = 18.93575989485557

$ws.Range("B4").Value = 5.507272574952419
$ws.Range("D4").Value = 8.226798883426534
$ws.Range("E4").Value = 13.22683167574835
$ws.Range("F4").Value = 37.0525373902735
$ws.Range("G4").Value = 41.84171940910331
$ws.Range("H4").Value = 17.46087381353415
$ws.Range("I4").Value = 20.52720030827601
$ws.Range("J4").Value = 10.08485727369825
$ws.Range("K4").Value = 14.57194764125702
$ws.Range("N4").Value = 18.97762182878245

$ws.Range("B5").Value = 5.479874969264014
$ws.Range("D5").Value = 8.219272242430566
$ws.Range("E5").Value = 13.21654669785974
$ws.Range("F5").Value = 37.04870973957625
$ws.Range("G5").Value = 41.81996873482925
$ws.Range("H5").Value = 17.47189424884663
$ws.Range("I5").Value = 20.55357519135529
$ws.Range("J5").Value = 10.08782104089962
$ws.Range("K5").Value = 14.48869610037789
$ws.Range("N5").Value = 18.99513952327376

$ws.Range("B6").Value = 5.475326282209581
$ws.Range("D6").Value = 8.218041664787084
$ws.Range("E6").Value = 13.21488106517028
$ws.Range("F6").Value = 37.04822915923859
$ws.Range("G6").Value = 41.81656993103298
$ws.Range("H6").Value = 17.47377886860014
$ws.Range("I6").Value = 20.55800361620473
$ws.Range("J6").Value = 10.0883391451227
$ws.Range("K6").Value = 14.47487575456184
$ws.Range("N6").Value = 18.99807606742038

$ws.Range("B7").Value = 5.506903061341346
$ws.Range("D7").Value = 8.226696091582978
$ws.Range("E7").Value = 13.22669014564817
$ws.Range("F7").Value = 37.05247537831821
$ws.Range("G7").Value = 41.84141180605574
$ws.Range("H7").Value = 17.46101877156688
$ws.Range("I7").Value = 20.52755273224425
$ws.Range("J7").Value = 10.08489550282986
$ws.Range("K7").Value = 14.57082472699112
$ws.Range("N7").Value = 18.97785621960525

$ws.Range("B8").Value = 5.646076002329227
$ws.Range("D8").Value = 8.267642923506267
$ws.Range("E8").Value = 13.2848818085804
$ws.Range("F8").Value = 37.09373787276891
$ws.Range("G8").Value = 41.98190618429698
$ws.Range("H8").Value = 17.41267548972939
$ws.Range("I8").Value = 20.40027683290328
$ws.Range("J8").Value = 10.07347534959448
$ws.Range("K8").Value = 14.99387396143098
$ws.Range("N8").Value = 18.89265263021509

$ws.Range("B9").Value = 5.91594504241557
$ws.Range("D9").Value = 8.357477621116194
$ws.Range("E9").Value = 13.4203125374924
$ws.Range("F9").Value = 37.25539894736926
$ws.Range("G9").Value = 42.36733683780901
$ws.Range("H9").Value = 17.34663660016767
$ws.Range("I9").Value = 20.17609965541263
$ws.Range("J9").Value = 10.06464980082853
$ws.Range("K9").Value = 15.81430378368819
$ws.Range("N9").Value = 18.7399205681546

$ws.Range("B10").Value = 6.110287215314509
$ws.Range("D10").Value = 8.428936928284569
$ws.Range("E10").Value = 13.53226306894778
$ws.Range("F10").Value = 37.42302171568416
$ws.Range("G10").Value = 42.71617956763777
$ws.Range("H10").Value = 17.31576619978287
$ws.Range("I10").Value = 20.02679364886821
$ws.Range("J10").Value = 10.06642922879455
$ws.Range("K10").Value = 16.40485201988621
$ws.Range("N10").Value = 18.63635660974951

$ws.Range("B11").Value = 6.19743676124112
$ws.Range("D11").Value = 8.462548628103788
$ws.Range("E11").Value = 13.58576424427437
$ws.Range("F11").Value = 37.50975718487668
$ws.Range("G11").Value = 42.88875226818151
$ws.Range("H11").Value = 17.30557513893325
$ws.Range("I11").Value = 19.96219422534584
$ws.Range("J11").Value = 10.06902866499
$ws.Range("K11").Value = 16.66956432868316
$ws.Range("N11").Value = 18.59109896139172

$ws.Range("B12").Value = 6.230225509238973
$ws.Range("D12").Value = 8.475427975701935
$ws.Range("E12").Value = 13.60638179529112
$ws.Range("F12").Value = 37.54409447838453
$ws.Range("G12").Value = 42.95605810749375
$ws.Range("H12").Value = 17.30227126196829
$ws.Range("I12").Value = 19.93820824504147
$ws.Range("J12").Value = 10.07026975836167
$ws.Range("K12").Value = 16.76913991145685
$ws.Range("N12").Value = 18.5742260764844

$ws.Range("B13").Value = 6.223173824145987
$ws.Range("D13").Value = 8.472647569294216
$ws.Range("E13").Value = 13.60192572603251
$ws.Range("F13").Value = 37.53663321338819
$ws.Range("G13").Value = 42.94147626399247
$ws.Range("H13").Value = 17.30295809691356
$ws.Range("I13").Value = 19.94335288301977
$ws.Range("J13").Value = 10.06999105984407
$ws.Range("K13").Value = 16.74772563314036
$ws.Range("N13").Value = 18.57784818153788

$ws.Range("B14").Value = 6.200138754751534
$ws.Range("D14").Value = 8.463605230739299
$ws.Range("E14").Value = 13.5874533643788
$ws.Range("F14").Value = 37.51255231763017
$ws.Range("G14").Value = 42.89425059653287
$ws.Range("H14").Value = 17.30529219196842
$ws.Range("I14").Value = 19.96021134222546
$ws.Range("J14").Value = 10.06912562920284
$ws.Range("K14").Value = 16.67777035219985
$ws.Range("N14").Value = 18.58970551197983

$ws.Range("B15").Value = 6.18600046604038
$ws.Range("D15").Value = 8.45808601796938
$ws.Range("E15").Value = 13.57863484604994
$ws.Range("F15").Value = 37.49799595143494
$ws.Range("G15").Value = 42.86557702284161
$ws.Range("H15").Value = 17.30679423919596
$ws.Range("I15").Value = 19.97059964763902
$ws.Range("J15").Value = 10.0686289417635
$ws.Range("K15").Value = 16.63483124061733
$ws.Range("N15").Value = 18.59700296483711

$ws.Range("B16").Value = 6.104563697750367
$ws.Range("D16").Value = 8.426761914265862
$ws.Range("E16").Value = 13.5288173999359
$ws.Range("F16").Value = 37.41756306455878
$ws.Range("G16").Value = 42.7051775642058
$ws.Range("H16").Value = 17.3165098420188
$ws.Range("I16").Value = 20.03108209921493
$ws.Range("J16").Value = 10.06629531572322
$ws.Range("K16").Value = 16.38746464727595
$ws.Range("N16").Value = 18.63935149184713

$ws.Range("B17").Value = 6.054259362141436
$ws.Range("D17").Value = 8.407823219449392
$ws.Range("E17").Value = 13.49890675131394
$ws.Range("F17").Value = 37.37089510791662
$ws.Range("G17").Value = 42.61030620959559
$ws.Range("H17").Value = 17.32345767331143
$ws.Range("I17").Value = 20.06903583433056
$ws.Range("J17").Value = 10.06532174588737
$ws.Range("K17").Value = 16.23463316700926
$ws.Range("N17").Value = 18.6658048149446

$ws.Range("B18").Value = 6.025209058628915
$ws.Range("D18").Value = 8.397034608548275
$ws.Range("E18").Value = 13.4819460248451
$ws.Range("F18").Value = 37.34504068986868
$ws.Range("G18").Value = 42.55704805933736
$ws.Range("H18").Value = 17.32781638827431
$ws.Range("I18").Value = 20.09117843443791
$ws.Range("J18").Value = 10.0649303021702
$ws.Range("K18").Value = 16.14636408619326
$ws.Range("N18").Value = 18.68119466877394

$ws.Range("B19").Value = 6.015354109366064
$ws.Range("D19").Value = 8.393399930119532
$ws.Range("E19").Value = 13.47624553156134
$ws.Range("F19").Value = 37.33645689578553
$ws.Range("G19").Value = 42.53924179475518
$ws.Range("H19").Value = 17.32935439062394
$ws.Range("I19").Value = 20.09872926402508
$ws.Range("J19").Value = 10.06482672843634
$ws.Range("K19").Value = 16.11641825255119
$ws.Range("N19").Value = 18.68643544023048

$ws.Range("B20").Value = 6.059626648250355
$ws.Range("D20").Value = 8.409828521802829
$ws.Range("E20").Value = 13.50206571994016
$ws.Range("F20").Value = 37.37576086378807
$ws.Range("G20").Value = 42.62027019122283
$ws.Range("H20").Value = 17.32268053631807
$ws.Range("I20").Value = 20.06496324784283
$ws.Range("J20").Value = 10.06540794632467
$ws.Range("K20").Value = 16.25094080421847
$ws.Range("N20").Value = 18.66297075314229

$ws.Range("B21").Value = 6.206910736814854
$ws.Range("D21").Value = 8.466257138701588
$ws.Range("E21").Value = 13.59169463828253
$ws.Range("F21").Value = 37.51958508640597
$ws.Range("G21").Value = 42.9080691653971
$ws.Range("H21").Value = 17.30459153263255
$ws.Range("I21").Value = 19.95524668265752
$ws.Range("J21").Value = 10.06937286511176
$ws.Range("K21").Value = 16.69833673075084
$ws.Range("N21").Value = 18.5862155402638

$ws.Range("B22").Value = 6.301914984019795
$ws.Range("D22").Value = 8.504014914835325
$ws.Range("E22").Value = 13.65235180260154
$ws.Range("F22").Value = 37.62227315095608
$ws.Range("G22").Value = 43.10754385690409
$ws.Range("H22").Value = 17.29600624169472
$ws.Range("I22").Value = 19.88631702610348
$ws.Range("J22").Value = 10.0734602634273
$ws.Range("K22").Value = 16.98681685270799
$ws.Range("N22").Value = 18.53759673537901

$ws.Range("B23").Value = 6.251334300874325
$ws.Range("D23").Value = 8.483785039104204
$ws.Range("E23").Value = 13.61979189679499
$ws.Range("F23").Value = 37.56667694214088
$ws.Range("G23").Value = 43.00005332058115
$ws.Range("H23").Value = 17.30029181208788
$ws.Range("I23").Value = 19.92285238264845
$ws.Range("J23").Value = 10.07114210091245
$ws.Range("K23").Value = 16.83323942829025
$ws.Range("N23").Value = 18.56340458909059

$ws.Range("B24").Value = 6.057200499654537
$ws.Range("D24").Value = 8.408921614134101
$ws.Range("E24").Value = 13.50063681636595
$ws.Range("F24").Value = 37.37355801534623
$ws.Range("G24").Value = 42.61576147036956
$ws.Range("H24").Value = 17.32303074515704
$ws.Range("I24").Value = 20.06680345953027
$ws.Range("J24").Value = 10.06536845089728
$ws.Range("K24").Value = 16.24356937278792
$ws.Range("N24").Value = 18.66425146701357

$ws.Range("B25").Value = 5.843484340108588
$ws.Range("D25").Value = 8.332189463960113
$ws.Range("E25").Value = 13.38144618941377
$ws.Range("F25").Value = 37.20305037199115
$ws.Range("G25").Value = 42.25141984664489
$ws.Range("H25").Value = 17.36141010510682
$ws.Range("I25").Value = 20.23403499717663
$ws.Range("J25").Value = 10.0655850263689
$ws.Range("K25").Value = 15.5940434929981
$ws.Range("N25").Value = 18.77971264600611
